$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks up front; we will rebuild them cleanly at the end
# so the relationship ids stay in a clean, predictable F2..F20 top-to-bottom order.
$ws.Hyperlinks.Delete()

# --- Insert blank rows for brand-new listings (ascending row order) ---
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(11).Insert()
$ws.Rows.Item(12).Insert()
$ws.Rows.Item(13).Insert()
$ws.Rows.Item(19).Insert()
$ws.Rows.Item(20).Insert()

# --- Row 2..20 field updates ---
$ws.Cells.Item(2,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(3,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(4,1).Value = "2025-10-01 12:38:31"
# Row 5: new listing
$ws.Cells.Item(5,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(5,2).Value = "【WordPressプラグイン開発】介護施設検索サイトの検索履歴等をMySQLに連携する開発者募集"
$ws.Cells.Item(5,3).Value = "システム開発"
$ws.Cells.Item(5,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(5,5).Value = "期限情報なし"
$ws.Cells.Item(5,6).Value = "https://www.lancers.jp/work/detail/5404680"
$ws.Cells.Item(5,7).Value = 133
$ws.Cells.Item(5,8).Value = "◆開発 ◇MySQL ○WordPress"

$ws.Cells.Item(6,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(7,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(8,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(9,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(10,1).Value = "2025-10-01 12:38:31"
# Row 11: new listing
$ws.Cells.Item(11,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(11,2).Value = "【急募】Salesforce・MA・CRMコンサルタント経験者を探しています!"
$ws.Cells.Item(11,3).Value = "システム開発"
$ws.Cells.Item(11,4).Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Cells.Item(11,5).Value = "期限情報なし"
$ws.Cells.Item(11,6).Value = "https://www.lancers.jp/work/detail/5371747"
$ws.Cells.Item(11,7).Value = 48
$ws.Cells.Item(11,8).Value = "◆コンサル"

# Row 12: new listing
$ws.Cells.Item(12,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(12,2).Value = "初回 Synology NAS DS925+ 導入・データ移行・アクセス制御再設計・5年間保守"
$ws.Cells.Item(12,3).Value = "システム開発"
$ws.Cells.Item(12,4).Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Cells.Item(12,5).Value = "期限情報なし"
$ws.Cells.Item(12,6).Value = "https://www.lancers.jp/work/detail/5404650"
$ws.Cells.Item(12,7).Value = 25

# Row 13: new listing
$ws.Cells.Item(13,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(13,2).Value = "〖リモート可〗Delphiエンジニア募集"
$ws.Cells.Item(13,3).Value = "システム開発"
$ws.Cells.Item(13,4).Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Cells.Item(13,5).Value = "期限情報なし"
$ws.Cells.Item(13,6).Value = "https://www.lancers.jp/work/detail/5341051"
$ws.Cells.Item(13,7).Value = 25

$ws.Cells.Item(14,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(15,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(16,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(17,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(18,1).Value = "2025-10-01 12:38:31"
# Row 19: new listing
$ws.Cells.Item(19,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(19,2).Value = "【急募】B型福祉施設の弁当集計表作成をお手伝いください!"
$ws.Cells.Item(19,3).Value = "システム開発"
$ws.Cells.Item(19,4).Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Cells.Item(19,5).Value = "期限情報なし"
$ws.Cells.Item(19,6).Value = "https://www.lancers.jp/work/detail/5404730"
$ws.Cells.Item(19,7).Value = 10

# Row 20: new listing
$ws.Cells.Item(20,1).Value = "2025-10-01 12:38:31"
$ws.Cells.Item(20,2).Value = "Meta Business Suite/Business設定の初期構築をサポート頂ける方を募集します"
$ws.Cells.Item(20,3).Value = "システム開発"
$ws.Cells.Item(20,4).Value = "~ 5,000 円 / 固定"
$ws.Cells.Item(20,5).Value = "期限情報なし"
$ws.Cells.Item(20,6).Value = "https://www.lancers.jp/work/detail/5404652"
$ws.Cells.Item(20,7).Value = 10


# --- Rebuild hyperlinks for column F, rows 2..20, in order ---
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5391872") | Out-Null
$ws.Range("F2").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5404026") | Out-Null
$ws.Range("F3").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5398112") | Out-Null
$ws.Range("F4").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5404680") | Out-Null
$ws.Range("F5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5404305") | Out-Null
$ws.Range("F6").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5403988") | Out-Null
$ws.Range("F7").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F8"), "https://www.lancers.jp/work/detail/5404059") | Out-Null
$ws.Range("F8").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F9"), "https://www.lancers.jp/work/detail/5404342") | Out-Null
$ws.Range("F9").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F10"), "https://www.lancers.jp/work/detail/5404426") | Out-Null
$ws.Range("F10").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F11"), "https://www.lancers.jp/work/detail/5371747") | Out-Null
$ws.Range("F11").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F12"), "https://www.lancers.jp/work/detail/5404650") | Out-Null
$ws.Range("F12").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F13"), "https://www.lancers.jp/work/detail/5341051") | Out-Null
$ws.Range("F13").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F14"), "https://www.lancers.jp/work/detail/5367840") | Out-Null
$ws.Range("F14").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F15"), "https://www.lancers.jp/work/detail/5372984") | Out-Null
$ws.Range("F15").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F16"), "https://www.lancers.jp/work/detail/5391221") | Out-Null
$ws.Range("F16").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F17"), "https://www.lancers.jp/work/detail/5404155") | Out-Null
$ws.Range("F17").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F18"), "https://www.lancers.jp/work/detail/5404010") | Out-Null
$ws.Range("F18").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F19"), "https://www.lancers.jp/work/detail/5404730") | Out-Null
$ws.Range("F19").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("F20"), "https://www.lancers.jp/work/detail/5404652") | Out-Null
$ws.Range("F20").Style = "Hyperlink"

# --- Column H width 13 -> 23 (account for the engine's ~0.8333 width offset) ---
$ws.Columns.Item(8).ColumnWidth = 22.166667
